# Updated cryptos list (refreshed Price / Volume(1h) snapshot), matching
# the upstream GitHub Actions scrape-and-commit job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '27.287.11'
$ws.Range('E2').Value = '  +0.45%  '

# Row 3 - Ethereum
$ws.Range('D3').Value = '1.776.31'
$ws.Range('E3').Value = '  +3.58%  '

# Row 4 - TetherUSD
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '

# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.00%  '

# Row 6 - USDC
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.06%  '

# Row 7 - XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5185'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +8.23%  '

# Row 8 - Cardano
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3688'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.70%  '

# Row 9 - OKB
$ws.Range('E9').Value = '  +1.61%  '

# Row 10 - Dogecoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07397'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.65%  '

# Row 11 - Polygon
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.088'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.99%  '

# Row 12 - BinanceUSD
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.06%  '

# Row 13 - Solana
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.48'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.05%  '

# Row 14 - Polkadot
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.072'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.51%  '

# Row 15 - WrappedEther
$ws.Range('D15').Value = '1.766.17'
$ws.Range('E15').Value = '  +3.15%  '

# Row 16 - Chainlink
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.968'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.32%  '

# Row 17 - Litecoin
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.38'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.48%  '

# Row 18 - ShibaInu
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001047'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.52%  '

# Row 19 - TRON
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06437'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.17%  '

# Row 20 - Dai
$ws.Range('E20').Value = '  +0.02%  '

# Row 21 - Avalanche
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.78'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.62%  '

# Row 22 - Uniswap
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.830'
$ws.Range('D22').Style = 'Normal'

# Row 23 - WrappedBTC
$ws.Range('D23').Value = '27.323.72'
$ws.Range('E23').Value = '  +0.45%  '

# Row 24 - Cosmos
$ws.Range('E24').Value = '  +3.81%  '

# Row 25 - Toncoin
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.121'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.51%  '

# Row 26 - Monero
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.04'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.88%  '

# Row 27 - EthereumClassic
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.98%  '

# Row 28 - LidoDAOToken
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.332'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +11.05%  '

# Row 29 - WrappedliquidstakedEther2.0
$ws.Range('D29').Value = '1.972.35'
$ws.Range('E29').Value = '  +3.26%  '

# Row 30 - BitcoinCash
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.05%  '

# Row 31 - ImmutableX
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.068'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.87%  '

# Row 32 - Stellar
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09784'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.79%  '

# Row 33 - Filecoin
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.582'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.74%  '

# Row 34 - HuobiToken
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.628'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.34%  '

# Row 35 - VeChain
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02246'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.94%  '

# Row 36 - Hedera
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05989'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.11%  '

# Row 37 - Aptos
$ws.Range('E37').Value = '  +1.39%  '

# Row 38 - InternetComputer(DFINITY)
$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6160'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.68%  '

# Row 39 - TheSandbox
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.846'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.88%  '

# Row 40 - Algorand
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2023'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.79%  '

# Row 41 - WEMIXTOKEN
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.434'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.15%  '

# Row 42 - FraxShare
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.086'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.09%  '

# Row 43 - TrustWalletToken
$ws.Range('E43').Value = '  +3.92%  '

# Row 44 - EnergySwap
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.10'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.74%  '

# Row 45 - Decentraland
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5772'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.58%  '

# Row 46 - PancakeSwap
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.632'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.49%  '

# Row 47 - Quant
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.27%  '

# Row 48 - NEARProtocol
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.889'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.58%  '

# Row 49 - EOS
$ws.Range('E49').Value = '  +2.75%  '

# Row 50 - Cronos
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06708'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.98%  '

# Row 51 - Aave
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.21%  '
